$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new schedule entry as row 30, continuing the table started in row 2.
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = 43639
$ws.Cells.Item(30, 3).Value = "The Strikers"
$ws.Cells.Item(30, 4).Value = "Westridge Warriors"
$ws.Cells.Item(30, 5).Value = "Garland Cricket Ground"
$ws.Cells.Item(30, 6).Value = "11.30 AM"
$ws.Cells.Item(30, 7).Value = "03.30 PM"

# Row 27 uses the same banding/style pattern the new row should follow, so
# copy its formatting down onto the freshly entered row.
$ws.Range("A27:G27").Copy()
$ws.Range("A30:G30").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Select the last entered cell and scroll the sheet so the new row is visible.
$ws.Range("D30").Select()

$wb.Save()
